$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strike through paragraph 2 ("Read writing R extensions...") including
#    the hyperlink run inside it.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Font.StrikeThrough = $true
$d.Hyperlinks.Item(1).Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 2) Strike through paragraph 4 ("Export only the functions...")
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 3) Strike through paragraph 5 ("Ask Rich whether...")
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------------
# 4) Insert "Worth removing dots..." as a brand-new paragraph right before
#    the bookmark paragraph (currently paragraph 7).
# ---------------------------------------------------------------------------
$rBefore = $d.Paragraphs.Item(7).Range
$rBefore.Collapse(1)
$rBefore.InsertBefore("Worth removing dots (“.”) where possible? `r")

# ---------------------------------------------------------------------------
# 5) Insert "Remove commented out lines" as text at the very start of the
#    (now-shifted) bookmark paragraph -- i.e. before the _GoBack bookmark.
# ---------------------------------------------------------------------------
$rStart = $d.Paragraphs.Item(8).Range
$rStart.Collapse(1)
$rStart.InsertBefore("Remove commented out lines")

# ---------------------------------------------------------------------------
# 6) Append the remaining new paragraphs (plain text first), each one
#    started via InsertAfter("`r...") from the very end of the document so
#    that no inherited character formatting can leak backward into earlier
#    paragraphs; we fix up character formatting afterward using explicit
#    Range(start, end) spans that exclude paragraph marks.
# ---------------------------------------------------------------------------
$endR = $d.Content
$endR.Collapse(0)
$endR.InsertAfter("`rUPLOAD")

$endR = $d.Content
$endR.Collapse(0)
$endR.InsertAfter("`rRun all updated versions on https://win-builder.r-project.org/upload.aspx in section “R-devel” – you’ll get an email with a link to all messages notes warnings etc. ")

$endR = $d.Content
$endR.Collapse(0)
$endR.InsertAfter("`rIn comments section, must state that package was archived on DATE, and that you have added S3 methods to namespace (look at Kurt Hornik email and other Dive emails to get times)")

$endR = $d.Content
$endR.Collapse(0)
$endR.InsertAfter("`r")

Write-Host "Paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" ("[" + $d.Paragraphs.Item($i).Range.Text + "]")
}

# ---------------------------------------------------------------------------
# 7) Make "UPLOAD" bold (text only, not the paragraph mark).
# ---------------------------------------------------------------------------
$pUpload = $d.Paragraphs.Item(9)
$uploadTextRange = $d.Range($pUpload.Range.Start, $pUpload.Range.End - 1)
Write-Host ("UPLOAD text range: [" + $uploadTextRange.Text + "]")
$uploadTextRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# 8) Turn the plain-text URL in the "Run all updated versions..." paragraph
#    into a real hyperlink.
# ---------------------------------------------------------------------------
$pRun = $d.Paragraphs.Item(10)
$urlText = "https://win-builder.r-project.org/upload.aspx"
$searchRange = $pRun.Range.Duplicate
$found = $searchRange.Find.Execute($urlText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found URL:" $found "range:" $searchRange.Start $searchRange.End "text:[" $searchRange.Text "]"
$d.Hyperlinks.Add($searchRange, $urlText, "", "", $urlText) | Out-Null
Write-Host "Hyperlinks count now:" $d.Hyperlinks.Count

